# Fixed a bug in catchSymbols
# Re-orders the symbol rows (rows 2 and 5-21) on the active sheet so that
# each row's data lands in the row position produced by the fixed
# catchSymbols logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @(301, 6, 45, 30, 60, 45)
    5  = @(101, 9, 30, 15, 60, 15)
    6  = @(902, 1, 0, 0, 0, 0)
    7  = @(1202, 2, 10, 10, 10, 10)
    8  = @(1203, 3, 15, 15, 15, 15)
    9  = @(1001, 18, 30, 75, 60, 72)
    10 = @(701, 3, 90, 45, 97, 15)
    11 = @(801, 3, 67, 65, 52, 45)
    12 = @(1201, 2, 10, 10, 10, 10)
    13 = @(901, 16, 15, 45, 60, 60)
    14 = @(601, 9, 60, 67, 60, 42)
    15 = @(401, 9, 48, 67, 75, 45)
    16 = @(502, 0, 4, 0, 0, 0)
    17 = @(1, 0, 2, 2, 2, 2)
    18 = @(3, 0, 3, 3, 3, 3)
    19 = @(2, 0, 2, 2, 2, 2)
    20 = @(1101, 0, 15, 30, 30, 0)
    21 = @(802, 0, 4, 5, 4, 0)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
}
